# تعديل تلقائي في شيت Card15 by admin at 2025-11-02 08:02:37
# Updates column A (card number) for rows 3 through 12 on the "Card15"
# worksheet from "2" to "15" to match the card identifier used in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

for ($row = 3; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = "15"
}
